$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.140.48"
$ws.Range("E2").Value = "  +3.48%  "

$ws.Range("D3").Value = "3.724.93"
$ws.Range("E3").Value = "  +6.69%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "421.13"
$ws.Range("E5").Value = "  +0.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.69"
$ws.Range("E6").Value = "  +0.56%  "

$ws.Range("D7").Value = "3.717.12"
$ws.Range("E7").Value = "  +6.74%  "

$ws.Range("E8").Value = "  +0.71%  "

$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.775"
$ws.Range("E10").Value = "  -0.30%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.185"
$ws.Range("E11").Value = "  +15.21%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000414"
$ws.Range("E12").Value = "  +62.42%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "43.07"
$ws.Range("E13").Value = "  -0.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.66"
$ws.Range("E14").Value = "  +7.53%  "

$ws.Range("D15").Value = "4.299.34"
$ws.Range("E15").Value = "  +5.74%  "

$ws.Range("E16").Value = "  -0.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.81"
$ws.Range("E17").Value = "  +2.31%  "

$ws.Range("D18").Value = "3.725.24"
$ws.Range("E18").Value = "  +6.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.42"
$ws.Range("E19").Value = "  +9.27%  "

$ws.Range("E20").Value = "  +4.44%  "

$ws.Range("D21").Value = "67.155.97"
$ws.Range("E21").Value = "  +3.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "447.05"
$ws.Range("E22").Value = "  -3.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.41"
$ws.Range("E23").Value = "  +24.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "90.01"
$ws.Range("E24").Value = "  +0.27%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.17"
$ws.Range("E25").Value = "  -0.92%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "37.85"
$ws.Range("E26").Value = "  +12.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.22"
$ws.Range("E27").Value = "  +3.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.34"
$ws.Range("E28").Value = "  -3.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.05"
$ws.Range("E29").Value = "  +4.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.75"
$ws.Range("E30").Value = "  +2.95%  "

$ws.Range("E31").Value = "  +9.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.77"
$ws.Range("E32").Value = "  -1.03%  "

$ws.Range("E33").Value = "  -2.94%  "

$ws.Range("E34").Value = "  +2.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "42.11"
$ws.Range("E35").Value = "  +5.73%  "

$ws.Range("E36").Value = "  +0.18%  "

$ws.Range("E37").Value = "  +0.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0498"
$ws.Range("E38").Value = "  +0.50%  "

$ws.Range("D39").Value = "0.0₃0758"
$ws.Range("E39").Value = "  +10.73%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.08"
$ws.Range("E40").Value = "  +32.66%  "

$ws.Range("E41").Value = "  +3.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "29.75"
$ws.Range("E42").Value = "  +37.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.997"
$ws.Range("E43").Value = "  -0.43%  "

$ws.Range("E44").Value = "  +5.18%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.26"
$ws.Range("E45").Value = "  +33.45%  "

$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "149.05"
$ws.Range("E46").Value = "  +2.43%  "

$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.13"
$ws.Range("E47").Value = "  +6.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.68"
$ws.Range("E48").Value = "  -2.59%  "

$ws.Range("E49").Value = "  -6.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.38"
$ws.Range("E50").Value = "  -1.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.311"
$ws.Range("E51").Value = "  -0.87%  "
